$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the crypto price/volume refresh.
# D-column price cells that look like plain numbers must keep their
# original text storage (t="inlineStr"/shared-string), so we pin the
# NumberFormat to text ("@") immediately before assigning those values
# to stop Excel auto-coercing them into floating point numbers.

$ws.Range("D2").Value = "30.149.67"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.921.90"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.47"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5079"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4077"
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08340"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.01"
$ws.Range("E11").Value = "  +5.04%  "
$ws.Range("D12").Value = "1.924.08"
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.432"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.256"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.004"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.72"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06516"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.55"
$ws.Range("E19").Value = "  +3.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").Value = "30.166.79"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.36"
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.195"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "2.143.15"
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("E26").Value = "  +4.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.77"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.263"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "128.93"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.134"
$ws.Range("E30").Value = "  +6.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1047"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.966"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.800"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.02452"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.322"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06438"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.216"
$ws.Range("E37").Value = "  +3.90%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2149"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6512"
$ws.Range("E39").Value = "  +3.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.590"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.211"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.45"
$ws.Range("E43").Value = "  +4.26%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.185"
$ws.Range("E44").Value = "  +9.34%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6054"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.625"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.77"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.211"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.135"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.93"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06835"
$ws.Range("E51").Value = "  +1.17%  "
